$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 20:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1604843
$ws.Range("C4").Value = 12120
$ws.Range("E4").Value = 1137392
$ws.Range("G4").Value = 720
$ws.Range("H4").Value = 95656

# Row 7 - España
$ws.Range("B7").Value = 280117
$ws.Range("C7").Value = 593
$ws.Range("E7").Value = 55219
$ws.Range("G7").Value = 52
$ws.Range("H7").Value = 27940

# Row 10 - Francia
$ws.Range("D10").Value = 63858
$ws.Range("E10").Value = 89502
$ws.Range("G10").Value = 83
$ws.Range("H10").Value = 28215

# Row 11 - Alemania
$ws.Range("B11").Value = 178876
$ws.Range("C11").Value = 345
$ws.Range("E11").Value = 12597
$ws.Range("G11").Value = 9
$ws.Range("H11").Value = 8279

# Row 25 - Ecuador
$ws.Range("B25").Value = 35306
$ws.Range("C25").Value = 452
$ws.Range("E25").Value = 28810
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = 2939

# Row 79 - Senegal
$ws.Range("E79").Value = 1529
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 32

# Row 111 - Republica de Chipre
$ws.Range("D111").Value = 561
$ws.Range("E111").Value = 345

# Row 188 - Botsuana
$ws.Range("D188").Value = 19
$ws.Range("E188").Value = 9
